# Updated cryptos list on Tue Mar 12 04:36:25 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and rotates the VeChain / FirstDigitalUSD / WEMIXToken rows (42-44).
#
# Numeric-looking text values (e.g. "541.97") are written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cell type) instead of auto-converting them to numbers; ClearFormats()
# then strips the transient "text" cell format Excel applies for that,
# so no stray style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.916.34'
$ws.Range('E2').Value = '  +4.72%  '
$ws.Range('D3').Value = '4.040.87'
$ws.Range('E3').Value = '  +4.56%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''541.97'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.51%  '
$ws.Range('D6').Value = '''152.63'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.18%  '
$ws.Range('D7').Value = '''0.693'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +13.95%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '''0.761'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.81%  '
$ws.Range('D10').Value = '''0.174'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.57%  '
$ws.Range('D11').Value = '''0.0000331'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.09%  '
$ws.Range('D12').Value = '''47.97'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +15.07%  '
$ws.Range('D13').Value = '''10.89'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.65%  '
$ws.Range('D14').Value = '4.694.00'
$ws.Range('E14').Value = '  +4.87%  '
$ws.Range('D15').Value = '4.017.88'
$ws.Range('E15').Value = '  +4.35%  '
$ws.Range('D16').Value = '''14.37'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('D17').Value = '''20.74'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.49%  '
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '71.900.97'
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('D21').Value = '''435.10'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.43%  '
$ws.Range('D22').Value = '''99.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +14.09%  '
$ws.Range('D23').Value = '''3.61'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').Value = '''4.29'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.95%  '
$ws.Range('D25').Value = '''14.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.32%  '
$ws.Range('D26').Value = '''11.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('D27').Value = '''11.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.38%  '
$ws.Range('D28').Value = '''37.12'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.14%  '
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('D30').Value = '''3.52'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +26.20%  '
$ws.Range('D31').Value = '''13.65'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.08%  '
$ws.Range('E32').Value = '  +5.52%  '
$ws.Range('D33').Value = '''687.87'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').Value = '''6.89'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').Value = '''67.18'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').Value = '''43.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +8.55%  '
$ws.Range('D37').Value = '''0.439'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('E38').Value = '  +6.65%  '
$ws.Range('D39').Value = '0.0₃0843'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').Value = '''3.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '''3.32'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.29%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0497'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.41%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''0.999'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '''0.152'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.75%  '
$ws.Range('D46').Value = '''2.72'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').Value = '''3.44'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('E48').Value = '  +10.35%  '
$ws.Range('D49').Value = '''3.05'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.40%  '
$ws.Range('D50').Value = '''3.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').Value = '''0.000272'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.79%  '
